# UserDomain.xlsx automation refresh: "User domain all scripts are working"
$wb = $excel.ActiveWorkbook

# --- addListItem sheet ---
$ws = $wb.Worksheets.Item("addListItem")
$ws.Range("D2").Value = "ADLILC.53186"
$ws.Range("A2").Value = "UserelevenPI"

# --- createUser sheet ---
$ws = $wb.Worksheets.Item("createUser")
$ws.Range("A2").Value = 120

# --- setHpDiary sheet ---
$ws = $wb.Worksheets.Item("setHpDiary")
$ws.Range("D2").Value = "30/10/2025"

# --- setHpClinicDiary sheet ---
$ws = $wb.Worksheets.Item("setHpClinicDiary")
$ws.Range("B2").Value = "28/10/2025"

# --- setHpLeave sheet ---
$ws = $wb.Worksheets.Item("setHpLeave")
$ws.Range("B2").Value = "25/10/2025"

# --- editHpLeave sheet ---
$ws = $wb.Worksheets.Item("editHpLeave")
$ws.Range("B2").Value = "13/10/2026"

# --- selection / view state updates ---
$ws = $wb.Worksheets.Item("setHpDiary")
$ws.Range("D6").Select()

$ws = $wb.Worksheets.Item("setHpLeave")
$ws.Range("B4").Select()

# --- active sheet switches from addListItem to createUser ---
$ws = $wb.Worksheets.Item("createUser")
$ws.Activate()
